$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "measure_id"
$ws.Range("B1").Value = "measure_qualifier"
$ws.Range("C1").Value = "Metric"
$ws.Range("D1").Value = "Description"
$ws.Range("E1").Value = "type"
$ws.Range("F1").Value = "achilles_note"
$ws.Range("G1").Value = "is_in_achilles"
$ws.Range("H1").Value = "achilles_analysis_id"

$ws.Range("A2").Value = "n"
$ws.Range("B2").Value = "n"
$ws.Range("C2").Value = "n"
$ws.Range("D2").Value = "Participant count"
$ws.Range("E2").Value = "simple"
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 1

$ws.Range("A3").Value = "earliest"
$ws.Range("B3").Value = "earliest"
$ws.Range("C3").Value = "earliest"
$ws.Range("D3").Value = "Earliest visit date in study"
$ws.Range("E3").Value = "simple"
$ws.Range("G3").Value = 0

$ws.Range("A4").Value = "latest"
$ws.Range("B4").Value = "latest"
$ws.Range("C4").Value = "latest"
$ws.Range("D4").Value = "Latest visit date in study"
$ws.Range("E4").Value = "simple"
$ws.Range("G4").Value = 0

$ws.Range("A5").Value = "span"
$ws.Range("B5").Value = "span"
$ws.Range("C5").Value = "span"
$ws.Range("D5").Value = "Time span of study (days)"
$ws.Range("E5").Value = "simple"
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 105

$ws.Range("A6").Value = "age_first_obs"
$ws.Range("B6").Value = "median"
$ws.Range("C6").Value = "enroll_age_median"
$ws.Range("D6").Value = "Median age of participants at first visit (years)"
$ws.Range("E6").Value = "patient_computation_then_aggregation"
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 103

$ws.Range("A7").Value = "age_first_obs"
$ws.Range("B7").Value = "mean"
$ws.Range("C7").Value = "enroll_age_mean"
$ws.Range("D7").Value = "Average age of participants at first visit (years)"
$ws.Range("E7").Value = "patient_computation_then_aggregation"
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 103

$ws.Range("A8").Value = "age_first_obs"
$ws.Range("B8").Value = "max"
$ws.Range("C8").Value = "enroll_age_max"
$ws.Range("D8").Value = "Oldest age of participants at first visit (years)"
$ws.Range("E8").Value = "patient_computation_then_aggregation"
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 103

$ws.Range("A9").Value = "age_first_obs"
$ws.Range("B9").Value = "min"
$ws.Range("C9").Value = "enroll_age_min"
$ws.Range("D9").Value = "Youngest  age of participants at first visit (years)"
$ws.Range("E9").Value = "patient_computation_then_aggregation"
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 103

$ws.Range("A10").Value = "age_first_obs"
$ws.Range("B10").Value = "sd"
$ws.Range("C10").Value = "enroll_age_sd"
$ws.Range("D10").Value = "Standard Deviation of ages of participants at first visit"
$ws.Range("E10").Value = "patient_computation_then_aggregation"
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 103

$ws.Range("A11").Value = "age_last_obs"
$ws.Range("B11").Value = "median"
$ws.Range("C11").Value = "end_age_median"
$ws.Range("D11").Value = "Median age of participants at last visit (years)"
$ws.Range("E11").Value = "patient_computation_then_aggregation"
$ws.Range("G11").Value = 0

$ws.Range("A12").Value = "age_last_obs"
$ws.Range("B12").Value = "mean"
$ws.Range("C12").Value = "end_age_mean"
$ws.Range("D12").Value = "Average age of participants at last visit (years)"
$ws.Range("E12").Value = "patient_computation_then_aggregation"
$ws.Range("G12").Value = 0

$ws.Range("A13").Value = "age_last_obs"
$ws.Range("B13").Value = "max"
$ws.Range("C13").Value = "end_age_max"
$ws.Range("D13").Value = "Oldest age of participants at last visit (years)"
$ws.Range("E13").Value = "patient_computation_then_aggregation"
$ws.Range("G13").Value = 0

$ws.Range("A14").Value = "age_last_obs"
$ws.Range("B14").Value = "min"
$ws.Range("C14").Value = "end_age_min"
$ws.Range("D14").Value = "Youngest  age of participants at last visit (years)"
$ws.Range("E14").Value = "patient_computation_then_aggregation"
$ws.Range("G14").Value = 0

$ws.Range("A15").Value = "age_last_obs"
$ws.Range("B15").Value = "sd"
$ws.Range("C15").Value = "end_age_sd"
$ws.Range("D15").Value = "Standard Deviation of ages of participants at last visit"
$ws.Range("E15").Value = "patient_computation_then_aggregation"
$ws.Range("G15").Value = 0

$ws.Range("A16").Value = "span_fl"
$ws.Range("B16").Value = "median"
$ws.Range("C16").Value = "span_median"
$ws.Range("D16").Value = "Median of the span between first and last visit of each patient (days)"
$ws.Range("E16").Value = "patient_computation_then_aggregation"
$ws.Range("G16").Value = 0

$ws.Range("A17").Value = "span_fl"
$ws.Range("B17").Value = "mean"
$ws.Range("C17").Value = "span_mean"
$ws.Range("D17").Value = "Average  of the span between first and last visit of each patient (days)"
$ws.Range("E17").Value = "patient_computation_then_aggregation"
$ws.Range("G17").Value = 0

$ws.Range("A18").Value = "span_fl"
$ws.Range("B18").Value = "sd"
$ws.Range("C18").Value = "span_sd"
$ws.Range("D18").Value = "Standard Deviation of the span between first and last visit of each patient (days)"
$ws.Range("E18").Value = "patient_computation_then_aggregation"
$ws.Range("G18").Value = 0

$ws.Range("A19").Value = "span_fl"
$ws.Range("B19").Value = "min"
$ws.Range("C19").Value = "span_min"
$ws.Range("D19").Value = "Smallest time span between first and last vist of a participant (days)"
$ws.Range("E19").Value = "patient_computation_then_aggregation"
$ws.Range("G19").Value = 0

$ws.Range("A20").Value = "span_fl"
$ws.Range("B20").Value = "max"
$ws.Range("C20").Value = "span_max"
$ws.Range("D20").Value = "Largest time span between first and last vist of a participant (days)"
$ws.Range("E20").Value = "patient_computation_then_aggregation"
$ws.Range("G20").Value = 0

$ws.Range("A21").Value = "visit_count"
$ws.Range("B21").Value = "median"
$ws.Range("C21").Value = "visit_count_median"
$ws.Range("D21").Value = "Median of number of visits per participant"
$ws.Range("E21").Value = "patient_computation_then_aggregation"
$ws.Range("G21").Value = 0

$ws.Range("A22").Value = "visit_count"
$ws.Range("B22").Value = "mean"
$ws.Range("C22").Value = "visit_count_mean"
$ws.Range("D22").Value = "Average of number of visits per participant"
$ws.Range("E22").Value = "patient_computation_then_aggregation"
$ws.Range("G22").Value = 0

$ws.Range("A23").Value = "visit_count"
$ws.Range("B23").Value = "sd"
$ws.Range("C23").Value = "visit_count_sd"
$ws.Range("D23").Value = "Standard Deviation of number of visits per participant"
$ws.Range("E23").Value = "patient_computation_then_aggregation"
$ws.Range("G23").Value = 0

$ws.Range("A24").Value = "visit_count"
$ws.Range("B24").Value = "min"
$ws.Range("C24").Value = "visit_count_min"
$ws.Range("D24").Value = ":east number of visits per participant"
$ws.Range("E24").Value = "patient_computation_then_aggregation"
$ws.Range("G24").Value = 0

$ws.Range("A25").Value = "visit_count"
$ws.Range("B25").Value = "max"
$ws.Range("C25").Value = "visit_count_max"
$ws.Range("D25").Value = "Most  visits per participant"
$ws.Range("E25").Value = "patient_computation_then_aggregation"
$ws.Range("G25").Value = 0

$ws.Range("E5").Select() | Out-Null